$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.213789594805377
$ws.Range("C2").Value = 0.07747796966823017
$ws.Range("D2").Value = 0.09677830275133914
$ws.Range("F2").Value = 2.503589715813845
$ws.Range("G2").Value = 1.788449332019695
$ws.Range("H2").Value = 1.570252067700395
$ws.Range("J2").Value = 0.2217237369403904
$ws.Range("K2").Value = 0.6879678976690116
$ws.Range("L2").Value = 0.3379112669186384
$ws.Range("N2").Value = 2.698121543173322
$ws.Range("B3").Value = 1.175233221520443
$ws.Range("C3").Value = 0.07546900259902145
$ws.Range("D3").Value = 0.09557695401960586
$ws.Range("F3").Value = 2.506502169072846
$ws.Range("G3").Value = 1.789873699694411
$ws.Range("H3").Value = 1.57598740455569
$ws.Range("J3").Value = 0.2225420270327803
$ws.Range("K3").Value = 0.6515379606852321
$ws.Range("L3").Value = 0.3337576971068188
$ws.Range("N3").Value = 2.719963048001844
$ws.Range("B4").Value = 1.152148396342426
$ws.Range("C4").Value = 0.07421664028713337
$ws.Range("D4").Value = 0.09487322769435735
$ws.Range("F4").Value = 2.509349661701179
$ws.Range("G4").Value = 1.791586846042222
$ws.Range("H4").Value = 1.580079002222845
$ws.Range("J4").Value = 0.2231139104880171
$ws.Range("K4").Value = 0.6294976425329537
$ws.Range("L4").Value = 0.3313529988256718
$ws.Range("N4").Value = 2.734114622693884
$ws.Range("B5").Value = 1.14288979695786
$ws.Range("C5").Value = 0.07370154841981247
$ws.Range("D5").Value = 0.09459502656767427
$ws.Range("F5").Value = 2.510776584822196
$ws.Range("G5").Value = 1.792495858416586
$ws.Range("H5").Value = 1.581889842491208
$ws.Range("J5").Value = 0.2233644484304449
$ws.Range("K5").Value = 0.6205988391944288
$ws.Range("L5").Value = 0.3304097872146912
$ws.Range("N5").Value = 2.740067796421076
$ws.Range("B6").Value = 1.141361408464576
$ws.Range("C6").Value = 0.07361573114922493
$ws.Range("D6").Value = 0.09454935065831194
$ws.Range("F6").Value = 2.511029628834095
$ws.Range("G6").Value = 1.792659537300352
$ws.Range("H6").Value = 1.582199200895644
$ws.Range("J6").Value = 0.2234071072625241
$ws.Range("K6").Value = 0.6191262119927785
$ws.Range("L6").Value = 0.3302553887640656
$ws.Range("N6").Value = 2.741067568353383
$ws.Range("B7").Value = 1.152022928994796
$ws.Range("C7").Value = 0.0742097127915855
$ws.Range("D7").Value = 0.0948694410026576
$ws.Range("F7").Value = 2.509367826194236
$ws.Range("G7").Value = 1.791598251381146
$ws.Range("H7").Value = 1.580102842729204
$ws.Range("J7").Value = 0.2231172184744885
$ws.Range("K7").Value = 0.6293772943861882
$ws.Range("L7").Value = 0.3313401295128457
$ws.Range("N7").Value = 2.734194154976706
$ws.Range("B8").Value = 1.200373442395573
$ws.Range("C8").Value = 0.07678918661768819
$ws.Range("D8").Value = 0.09635706795690169
$ws.Range("F8").Value = 2.504374196470593
$ws.Range("G8").Value = 1.788766416899975
$ws.Range("H8").Value = 1.572111368730788
$ws.Range("J8").Value = 0.2219914850497204
$ws.Range("K8").Value = 0.6753390300273452
$ws.Range("L8").Value = 0.3364489570099209
$ws.Range("N8").Value = 2.705498749003766
$ws.Range("B9").Value = 1.299843253204529
$ws.Range("C9").Value = 0.08169837629115762
$ws.Range("D9").Value = 0.09954149828890735
$ws.Range("F9").Value = 2.502980041364737
$ws.Range("G9").Value = 1.789867860352075
$ws.Range("H9").Value = 1.560958421888117
$ws.Range("J9").Value = 0.2203338622012829
$ws.Range("K9").Value = 0.7680597200325963
$ws.Range("L9").Value = 0.3476188407624932
$ws.Range("N9").Value = 2.655103857766044
$ws.Range("B10").Value = 1.375745685506843
$ws.Range("C10").Value = 0.08521507945158646
$ws.Range("D10").Value = 0.1020418581564186
$ws.Range("F10").Value = 2.507070656217394
$ws.Range("G10").Value = 1.794738099603421
$ws.Range("H10").Value = 1.555513188041459
$ws.Range("J10").Value = 0.2194498624778696
$ws.Range("K10").Value = 0.8377532890818316
$ws.Range("L10").Value = 0.3565234096262202
$ws.Range("N10").Value = 2.621658626450021
$ws.Range("B11").Value = 1.410885710874368
$ws.Range("C11").Value = 0.08679555742633482
$ws.Range("D11").Value = 0.1032138146360495
$ws.Range("F11").Value = 2.510041460190152
$ws.Range("G11").Value = 1.797836622028129
$ws.Range("H11").Value = 1.553631804114573
$ws.Range("J11").Value = 0.2191199154045762
$ws.Range("K11").Value = 0.8697989499856931
$ws.Range("L11").Value = 0.3607251577979014
$ws.Range("N11").Value = 2.607220311314286
$ws.Range("B12").Value = 1.424279848419644
$ws.Range("C12").Value = 0.08739127871940156
$ws.Range("D12").Value = 0.1036625268964073
$ws.Range("F12").Value = 2.511325917779558
$ws.Range("G12").Value = 1.799136973482533
$ws.Range("H12").Value = 1.553004929649859
$ws.Range("J12").Value = 0.219005329190864
$ws.Range("K12").Value = 0.8819826834128719
$ws.Range("L12").Value = 0.3623378775581614
$ws.Range("N12").Value = 2.601864500842524
$ws.Range("B13").Value = 1.421391304531937
$ws.Range("C13").Value = 0.08726310270507298
$ws.Range("D13").Value = 0.1035656707203216
$ws.Range("F13").Value = 2.511042196095559
$ws.Range("G13").Value = 1.798851270310095
$ws.Range("H13").Value = 1.553136134233895
$ws.Range("J13").Value = 0.219029547091786
$ws.Range("K13").Value = 0.8793565371879879
$ws.Range("L13").Value = 0.3619895901391601
$ws.Range("N13").Value = 2.603013002992007
$ws.Range("B14").Value = 1.411985906865596
$ws.Range("C14").Value = 0.0868446233065896
$ws.Range("D14").Value = 0.1032506321766249
$ws.Range("F14").Value = 2.510143937854508
$ws.Range("G14").Value = 1.797941057113306
$ws.Range("H14").Value = 1.553578516458387
$ws.Range("J14").Value = 0.2191102808781622
$ws.Range("K14").Value = 0.8708003385910956
$ws.Range("L14").Value = 0.3608574047594715
$ws.Range("N14").Value = 2.606777447220708
$ws.Range("B15").Value = 1.406236190550146
$ws.Range("C15").Value = 0.08658793173577806
$ws.Range("D15").Value = 0.1030583010004023
$ws.Range("F15").Value = 2.509614493657367
$ws.Range("G15").Value = 1.797400066387866
$ws.Range("H15").Value = 1.553860628661909
$ws.Range("J15").Value = 0.2191610808339064
$ws.Range("K15").Value = 0.8655657550492322
$ws.Range("L15").Value = 0.3601667194749325
$ws.Range("N15").Value = 2.609097822236187
$ws.Range("B16").Value = 1.373461460224689
$ws.Range("C16").Value = 0.08511140356143443
$ws.Range("D16").Value = 0.1019659588058559
$ws.Range("F16").Value = 2.506898831096564
$ws.Range("G16").Value = 1.794553375171873
$ws.Range("H16").Value = 1.555648119506714
$ws.Range("J16").Value = 0.2194728759422233
$ws.Range("K16").Value = 0.8356658715933065
$ws.Range("L16").Value = 0.3562518451701635
$ws.Range("N16").Value = 2.622617826704854
$ws.Range("B17").Value = 1.353511465188944
$ws.Range("C17").Value = 0.0842006627936982
$ws.Range("D17").Value = 0.1013046523285368
$ws.Range("F17").Value = 2.50551705320818
$ws.Range("G17").Value = 1.793033207597901
$ws.Range("H17").Value = 1.556897188412492
$ws.Range("J17").Value = 0.2196826249540251
$ws.Range("K17").Value = 0.8174104961321405
$ws.Range("L17").Value = 0.3538887990534505
$ws.Range("N17").Value = 2.63111073191655
$ws.Range("B18").Value = 1.342094359568364
$ws.Range("C18").Value = 0.08367501319278148
$ws.Range("D18").Value = 0.1009275393200397
$ws.Range("F18").Value = 2.504826773868217
$ws.Range("G18").Value = 1.792241960126745
$ws.Range("H18").Value = 1.557671693527865
$ws.Range("J18").Value = 0.2198100631621926
$ws.Range("K18").Value = 0.8069426739301662
$ws.Range("L18").Value = 0.3525438580609404
$ws.Range("N18").Value = 2.636068671306369
$ws.Range("B19").Value = 1.338238634954848
$ws.Range("C19").Value = 0.0834967253721004
$ws.Range("D19").Value = 0.1008004154922517
$ws.Range("F19").Value = 2.504611006236686
$ws.Range("G19").Value = 1.791988331121942
$ws.Range("H19").Value = 1.55794356123937
$ws.Range("J19").Value = 0.2198543794971748
$ws.Range("K19").Value = 0.803403989525691
$ws.Range("L19").Value = 0.3520909298178196
$ws.Range("N19").Value = 2.6377598880024
$ws.Range("B20").Value = 1.35562921863027
$ws.Range("C20").Value = 0.08429780059508118
$ws.Range("D20").Value = 0.1013747132056224
$ws.Range("F20").Value = 2.50565333321228
$ws.Range("G20").Value = 1.79318642995932
$ws.Range("H20").Value = 1.556758420151922
$ws.Range("J20").Value = 0.2196595936082844
$ws.Range("K20").Value = 0.8193504836875434
$ws.Range("L20").Value = 0.3541388783403221
$ws.Range("N20").Value = 2.630199087310004
$ws.Range("B21").Value = 1.414746133491008
$ws.Range("C21").Value = 0.08696761604753789
$ws.Range("D21").Value = 0.1033430335448031
$ws.Range("F21").Value = 2.510403450919412
$ws.Range("G21").Value = 1.798204961578847
$ws.Range("H21").Value = 1.55344625660787
$ws.Range("J21").Value = 0.2190862864889525
$ws.Range("K21").Value = 0.8733121819968801
$ws.Range("L21").Value = 0.3611893694408082
$ws.Range("N21").Value = 2.605668707572676
$ws.Range("B22").Value = 1.453891461263538
$ws.Range("C22").Value = 0.08869634185438713
$ws.Range("D22").Value = 0.1046580872543501
$ws.Range("F22").Value = 2.514437430429012
$ws.Range("G22").Value = 1.802225178287841
$ws.Range("H22").Value = 1.551780261115994
$ws.Range("J22").Value = 0.2187719584732477
$ws.Range("K22").Value = 0.9088631633772479
$ws.Range("L22").Value = 0.3659231768572653
$ws.Range("N22").Value = 2.590287594866517
$ws.Range("B23").Value = 1.432952479928019
$ws.Range("C23").Value = 0.08777516563688437
$ws.Range("D23").Value = 0.1039536132343102
$ws.Range("F23").Value = 2.512199413455761
$ws.Range("G23").Value = 1.800011763716043
$ws.Range("H23").Value = 1.552623834108914
$ws.Range("J23").Value = 0.2189342057673684
$ws.Range("K23").Value = 0.8898630913687953
$ws.Range("L23").Value = 0.3633851697724566
$ws.Range("N23").Value = 2.598437203457266
$ws.Range("B24").Value = 1.354671618996861
$ws.Range("C24").Value = 0.08425389098990621
$ws.Range("D24").Value = 0.101343029084056
$ws.Range("F24").Value = 2.505591396668777
$ws.Range("G24").Value = 1.793116900468306
$ws.Range("H24").Value = 1.55682098159059
$ws.Range("J24").Value = 0.2196699847385553
$ws.Range("K24").Value = 0.8184733298229787
$ws.Range("L24").Value = 0.3540257751074876
$ws.Range("N24").Value = 2.630611007414252
$ws.Range("B25").Value = 1.272437345636916
$ws.Range("C25").Value = 0.08038619035752959
$ws.Range("D25").Value = 0.09865164118435388
$ws.Range("F25").Value = 2.502458616691271
$ws.Range("G25").Value = 1.788857039828684
$ws.Range("H25").Value = 1.563492462878884
$ws.Range("J25").Value = 0.2207235640227729
$ws.Range("K25").Value = 0.7426997983484114
$ws.Range("L25").Value = 0.3444742231896214
$ws.Range("N25").Value = 2.668108162103607
